$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dbTesting")

# A2: client name test-data value bumped for a fresh test run (new shared string)
$ws.Range("A2").Value = "Gandhali16"

# Column A: size to fit the "clinet name" header (mirrors Sheet1's column A sizing)
$ws.Columns.Item(1).ColumnWidth = 10.6

# Move/collapse the view: selection on E10, no scrolled-in topLeftCell
$ws.Activate()
$ws.Range("E10").Select()
